# Acerto do individuo que anda pela curva
# Adds a new task row (26) to the "Planilha1" sheet:
#   A26 = "Corrigir bugs nos individuos"  (new shared string)
#   B26 = "Douglas"                       (same dev as the row above)
# The new row reuses the formatting of row 25 (style s="8" on column A,
# s="13" on column B), and the final selection ends up on C26 - mirroring
# what happens when a user fills A26/B26 in Excel and tabs to C26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row onto the new row
# before filling in values, so the new cells pick up the same styles
# (column A has no sheet-wide default style, so this must be explicit).
$ws.Range("A25").Copy()
$ws.Range("A26").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A26").Value = "Corrigir bugs nos individuos"
$ws.Range("B26").Value = "Douglas"

# Leave the selection where the user would land after entering the data.
$ws.Range("C26").Select()
